$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)

$ws.Range("D2").Value = "30.268.14"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").Value = "2.088.46"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4404"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.18"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09319"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.674"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.893"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.88%  "

$ws.Range("D15").Value = "2.062.64"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001155"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.005"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06676"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.355"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").Value = "30.283.27"
$ws.Range("E23").Value = "  +1.70%  "

$ws.Range("E24").Value = "  -0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.297"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.514"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.130"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.666"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.209"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.664"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.855"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02625"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06767"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6988"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6803"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.336"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.370"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.626"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.211"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.214"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
